$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "331.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.01%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "45.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.50%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.602"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.52%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08354"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.25%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.039"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.81%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9745"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.57%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.560"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.92%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1153"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.24%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1928"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.33%"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.64%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09944"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.56%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04693"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.76%"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.52%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001286"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.47%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006023"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.95%"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.37%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.450"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.80%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-3.18%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.94%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2651"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.10%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04199"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.17%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001311"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.71%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004584"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "5.78%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001303"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.54%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.05%"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "6.67%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05794"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.59%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007663"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.37%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1436"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.89%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007298"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.50%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002118"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "5.10%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008210"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.91%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3403"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007305"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.77%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.13%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0005809"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.05%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003507"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-7.20%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.003503"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.74%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00002103"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.13%"

